# Generate Report for Handback
#
# For the "344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.md" source file, the latest
# handoff round-tripped and came back unchanged ("in sync with en-US"), so:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Target File / Latest Handback File are filled in (same file that
#     was targeted / handed back)
#   - Latest Handback DateTime is stamped with the handback time
# This applies to both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = "Handed back: in sync with en-US"

$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d78f89cd38dfebc73e0cf0543c3f401d75d8dd20/e2e/344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.md",
    [Type]::Missing,
    [Type]::Missing,
    "344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.md"
)
$ws.Range("E2").Style = "HyperLink"
$ws.Range("E2").Font.Underline = 2
$ws.Range("E2").Font.Color = 15570276

$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fab7ba74105183b473e8ab209e26d28779cedc59/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.011397bac90f62833f853d530b22560bfee78cac.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.011397bac90f62833f853d530b22560bfee78cac.zh-cn.xlf"
)
$ws.Range("F2").Style = "HyperLink"
$ws.Range("F2").Font.Underline = 2
$ws.Range("F2").Font.Color = 15570276

$ws.Range("G2").Value = "2016-03-09 14:08:28"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = "Handed back: in sync with en-US"

$ws2.Hyperlinks.Add(
    $ws2.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d78f89cd38dfebc73e0cf0543c3f401d75d8dd20/e2e/344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.md",
    [Type]::Missing,
    [Type]::Missing,
    "344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.md"
)
$ws2.Range("E2").Style = "HyperLink"
$ws2.Range("E2").Font.Underline = 2
$ws2.Range("E2").Font.Color = 15570276

$ws2.Hyperlinks.Add(
    $ws2.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c74f7358e2358f5380516d13472cf2ad39bb1f7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.011397bac90f62833f853d530b22560bfee78cac.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "344e85f5-c7e4-4c00-b77c-eb673ed3bbf4.011397bac90f62833f853d530b22560bfee78cac.de-de.xlf"
)
$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("F2").Font.Underline = 2
$ws2.Range("F2").Font.Color = 15570276

$ws2.Range("G2").Value = "2016-03-09 14:08:33"

Write-Host "Handback report generated."
